$wb = $excel.ActiveWorkbook

# --- Update Status for the "5e71d6e3-...md" file (row 3) from
# "Ready for handoff" to "Handback transform failed" everywhere it appears:
# the Overview sheet (both language status columns) and each language
# sheet's own Status column.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- Populate the "Error Detail" column (K) for row 3 on each language
# sheet with the handback-transform failure detail.
$wsZhCn.Range("K3").Value = "Handback file name: rfvur4cd.j4h is different with handoff file name: 5e71d6e3-516f-4fe6-b5a8-cfba936fbab7.96118638f638cc8672fe55f8ae6878d4ec3a33c2.zh-cn."

$wsDeDe.Range("K3").Value = "Handback file name: rfvur4cd.j4h is different with handoff file name: 5e71d6e3-516f-4fe6-b5a8-cfba936fbab7.96118638f638cc8672fe55f8ae6878d4ec3a33c2.de-de."
